$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.463.80'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.05%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.992.94'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -4.86%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.015'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +1.22%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '329.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.96%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5005'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.59%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4233'
$ws.Range('D8').Style = 'Normal'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '54.35'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.21%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08920'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.76%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.111'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.05%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.28'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -6.09%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.005.43'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.25%  '

$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.962'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.30%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.455'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.72%  '

$ws.Range('E16').Value = '  +1.15%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '94.09'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.28%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001111'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.17%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06771'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.44%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -8.16%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.928'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -6.31%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '29.487.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.06%  '

$ws.Range('E24').Value = '  -3.84%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.319'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.39%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '20.77'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.96%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '156.89'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.71%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.279'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.59%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.303'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.28%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.74'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.27%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.059'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.08%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09951'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.17%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.545'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.96%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.834'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.95%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.802'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.99%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02454'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.81%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.199'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -9.42%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06378'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.50%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.296'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.43%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6539'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.65%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.60'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.99%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2039'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.10%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.011'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.98%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6333'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.42%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.59'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.63%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.205'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.97%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.311'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.27%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.499'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.75%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000340'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.69%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06955'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.05%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.132'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -8.01%  '
